$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 19: Unbreak My Heart
$ws.Range("H19").Value = 536.5789
$ws.Range("I19").Value = 589
$ws.Range("J19").Value = 489.4
$ws.Range("K19").Value = 589
$ws.Range("L19").Value = 489.4
$ws.Range("M19").Value = -414
$ws.Range("N19").Value = -839.4

# ALC!row 129: Practical Command
$ws.Range("H129").Value = 853.5862
$ws.Range("J129").Value = 868.38464
$ws.Range("L129").Value = 2605.15392
$ws.Range("N129").Value = -12605.15392

# ALC!row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2599.578
$ws.Range("I137").Value = 2180.2727
$ws.Range("J137").Value = 3000.652
$ws.Range("K137").Value = 6540.8181
$ws.Range("L137").Value = 9001.956
$ws.Range("M137").Value = -3990.8181
$ws.Range("N137").Value = -14101.956

# ALC!row 140: Tome for Tradition
$ws.Range("H140").Value = 37463.625
$ws.Range("J140").Value = 37463.625
$ws.Range("L140").Value = 37463.625
$ws.Range("N140").Value = -47823.625

# ALC!row 141: Remedy for Reason
$ws.Range("H141").Value = 1160
$ws.Range("I141").Value = 1160
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3480
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1700
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 32: Ingot We Trust
$ws.Range("H32").Value = 9836.138999999999
$ws.Range("I32").Value = 7136.273
$ws.Range("J32").Value = 22064.941
$ws.Range("K32").Value = 7136.273
$ws.Range("L32").Value = 22064.941
$ws.Range("M32").Value = -6849.273
$ws.Range("N32").Value = -22638.941

# ARM!row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 125001870
$ws.Range("I61").Value = 166668160
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 166668160
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -166667948
$ws.Range("N61").Value = -3424

# ARM!row 134: Brace for More Vambraces
$ws.Range("H134").Value = 36300
$ws.Range("J134").Value = 36300
$ws.Range("L134").Value = 36300
$ws.Range("N134").Value = -46440

# ARM!row 136: Metal with Mettle
$ws.Range("H136").Value = 125001870
$ws.Range("I136").Value = 166668160
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 500004480
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -500001930
$ws.Range("N136").Value = -14100

# ARM!row 140: A Hand for a Deckhand
$ws.Range("H140").Value = 35725.25
$ws.Range("J140").Value = 35725.25
$ws.Range("L140").Value = 35725.25
$ws.Range("N140").Value = -46085.25

# ARM!row 141: Essays on Equipment
$ws.Range("H141").Value = 41359.5
$ws.Range("J141").Value = 41359.5
$ws.Range("L141").Value = 41359.5
$ws.Range("N141").Value = -51719.5

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 81: Diamond Sawdust
$ws.Range("H81").Value = 19991
$ws.Range("J81").Value = 19991
$ws.Range("L81").Value = 19991
$ws.Range("N81").Value = -22113

# BSM!row 84: I'm a Lumberjack and I'm Okay (L)
$ws.Range("H84").Value = 19991
$ws.Range("J84").Value = 19991
$ws.Range("L84").Value = 59973
$ws.Range("N84").Value = -70581

# BSM!row 107: The Gold Experience
$ws.Range("H107").Value = 874.5
$ws.Range("I107").Value = 790.6923
$ws.Range("K107").Value = 790.6923
$ws.Range("M107").Value = 1129.3077

# BSM!row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4335.1333
$ws.Range("I134").Value = 1016.3077
$ws.Range("J134").Value = 25907.5
$ws.Range("K134").Value = 3048.9231
$ws.Range("L134").Value = 77722.5
$ws.Range("M134").Value = -513.9231
$ws.Range("N134").Value = -82792.5

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 51099.5
$ws.Range("I58").Value = 2200
$ws.Range("K58").Value = 2200
$ws.Range("M58").Value = -1997

# CRP!row 103: Spare a Rod and Spoil the Fishers
$ws.Range("H103").Value = 20960.572
$ws.Range("I103").Value = 16841.334
$ws.Range("J103").Value = 24050
$ws.Range("K103").Value = 16841.334
$ws.Range("L103").Value = 24050
$ws.Range("M103").Value = -15669.334
$ws.Range("N103").Value = -26394

# CRP!row 134: Wood You Be Quiet
$ws.Range("H134").Value = 20001630
$ws.Range("I134").Value = 1700
$ws.Range("K134").Value = 5100
$ws.Range("M134").Value = -2565

# CRP!row 136: Turali Quality
$ws.Range("H136").Value = 51099.5
$ws.Range("I136").Value = 2200
$ws.Range("K136").Value = 6600
$ws.Range("M136").Value = -4050

# CRP!row 140: Spear Pressure
$ws.Range("H140").Value = 56425
$ws.Range("J140").Value = 56425
$ws.Range("L140").Value = 56425
$ws.Range("N140").Value = -66785

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 11: Putting the Squeeze On
$ws.Range("H11").Value = 220
$ws.Range("I11").Value = 220
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 660
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -520
$ws.Range("N11").ClearContents()

# CUL!row 56: Culture Club
$ws.Range("H56").Value = 7159.909
$ws.Range("I56").Value = 7159.909
$ws.Range("K56").Value = 7159.909
$ws.Range("M56").Value = -6629.909

# CUL!row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 38636.25
$ws.Range("I140").Value = 86249.086
$ws.Range("J140").Value = 2926.625
$ws.Range("K140").Value = 258747.258
$ws.Range("L140").Value = 8779.875
$ws.Range("M140").Value = -253567.258
$ws.Range("N140").Value = -19139.875

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 141: Mask Maker
$ws.Range("H141").Value = 54279.875
$ws.Range("J141").Value = 54279.875
$ws.Range("L141").Value = 54279.875
$ws.Range("N141").Value = -64639.875

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 136: Respect for Br'aax
$ws.Range("H136").Value = 1490.05
$ws.Range("I136").Value = 1323.3889
$ws.Range("J136").Value = 2990
$ws.Range("K136").Value = 3970.1667
$ws.Range("L136").Value = 8970
$ws.Range("M136").Value = -1420.1667
$ws.Range("N136").Value = -14070

# LTW!row 140: Worqor Zormor or Bust
$ws.Range("H140").Value = 49806.668
$ws.Range("J140").Value = 49806.668
$ws.Range("L140").Value = 49806.668
$ws.Range("N140").Value = -60166.668

# LTW!row 141: Just Generally Freezing
$ws.Range("J141").Value = 46061.668
$ws.Range("L141").Value = 46061.668
$ws.Range("N141").Value = -56421.668

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 46: Crunching the Numbers
$ws.Range("H46").Value = 40999.2
$ws.Range("J46").Value = 40999.2
$ws.Range("L46").Value = 40999.2
$ws.Range("N46").Value = -41461.2

# WVR!row 99: Say Yes to Formal Dress
$ws.Range("H99").Value = 15500
$ws.Range("J99").Value = 15500
$ws.Range("L99").Value = 15500
$ws.Range("N99").Value = -21490

# WVR!row 100: Of Great Import
$ws.Range("H100").Value = 380.6
$ws.Range("I100").Value = 383.33334
$ws.Range("J100").Value = 376.5
$ws.Range("K100").Value = 766.66668
$ws.Range("L100").Value = 753
$ws.Range("M100").Value = -225.66668
$ws.Range("N100").Value = -1835

# WVR!row 134: Cloth for Canvas
$ws.Range("H134").Value = 40999.2
$ws.Range("J134").Value = 40999.2
$ws.Range("L134").Value = 122997.6
$ws.Range("N134").Value = -128067.6

# WVR!row 136: Weaving the Envelope
$ws.Range("H136").Value = 2020.5
$ws.Range("I136").Value = 1475
$ws.Range("J136").Value = 2384.1667
$ws.Range("K136").Value = 4425
$ws.Range("L136").Value = 7152.500100000001
$ws.Range("M136").Value = -1875
$ws.Range("N136").Value = -12252.5001

# WVR!row 140: Glamorous Gloves
$ws.Range("H140").Value = 35460
$ws.Range("J140").Value = 35460
$ws.Range("L140").Value = 35460
$ws.Range("N140").Value = -45820

# WVR!row 141: Silk for Sunperch
$ws.Range("H141").Value = 57003.332
$ws.Range("J141").Value = 57003.332
$ws.Range("L141").Value = 57003.332
$ws.Range("N141").Value = -67363.33199999999
